# Add a "Treatment" column as the new column B, shifting the existing
# analyte columns (old B..AH) one column to the right (new C..AI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; Excel shifts B:AH -> C:AI and extends
# the sheet dimension (A1:AH28 -> A1:AI28) automatically.
$ws.Columns("B:B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "Treatment"

# Per-sample treatment assignment (row 1 is the header row).
$treatments = @{
  2  = "untreated"
  3  = "actigard"
  4  = "actigard"
  5  = "actigard"
  6  = "actigard"
  7  = "actigard"
  8  = "untreated"
  9  = "untreated"
  10 = "untreated"
  11 = "untreated"
  12 = "untreated"
  13 = "untreated"
  14 = "untreated"
  15 = "untreated"
  16 = "untreated"
  17 = "untreated"
  18 = "rrv"
  19 = "rrv"
  20 = "rrv"
  21 = "rrv"
  22 = "rrv"
  23 = "rrv"
  24 = "rrv"
  25 = "rrv"
  26 = "rrv"
  27 = "rrv"
  28 = "rrv"
}

foreach ($row in $treatments.Keys) {
  $ws.Range("B$row").Value = $treatments[$row]
}
